# Update scripts with new TPM: recomputed NATMI edge-expression / specificity
# statistics, and the "Sending cluster" now reads ECs (was MuSCs) for rows 2-5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row = 2; Cells = @{ A = "ECs"; D = "ECs"; G = 0.1161373333333333; H = 0.348412; I = 0.0005423317774654872; J = 0.0005423317774654872; M = 3.078094333333333; N = 9.234283; O = 0.1162262227649079; P = 0.1162262227649079; Q = 0.3574816676217777; R = 3.217335008596; S = 0.000063033173980192168817798803; T = 0.00006303317398019218237032596 } },
    @{ Row = 3; Cells = @{ A = "ECs"; D = "FAPs"; G = 0.1161373333333333; H = 0.348412; I = 0.0005423317774654872; J = 0.0005423317774654872; O = 0.390990551238643; P = 0.390990551238643; Q = 1.202585362890667; R = 10.823268266016; S = 0.0002120466006254639; T = 0.0002120466006254639 } },
    @{ Row = 4; Cells = @{ A = "ECs"; D = "MuSCs"; G = 0.1161373333333333; H = 0.348412; I = 0.0005423317774654872; J = 0.0005423317774654872; M = 12.706793; N = 38.120379; O = 0.4797976910104138; P = 0.4797976910104139; Q = 1.475733054238667; R = 13.281597488148; S = 0.0002602095345895143; T = 0.0002602095345895144 } },
    @{ Row = 5; Cells = @{ A = "ECs"; G = 0.1161373333333333; H = 0.348412; I = 0.0005423317774654872; J = 0.0005423317774654872; M = 0.3439043333333334; N = 1.031713; O = 0.01298553498603535; P = 0.01298553498603535; Q = 0.03994013219511111; R = 0.3594611897560001; S = 0.000007042468270316821034923022; T = 0.000007042468270316821034923022 } },
    @{ Row = 6; Cells = @{ D = "ECs"; G = 214.0283; H = 642.0849000000001; I = 0.9994576682225345; J = 0.9994576682225345; M = 3.078094333333333; N = 9.234283; O = 0.1162262227649079; P = 0.1162262227649079; Q = 658.7992974029667; R = 5929.1936766267; S = 0.1161631895909277; T = 0.1161631895909277 } },
    @{ Row = 7; Cells = @{ D = "FAPs"; G = 214.0283; H = 642.0849000000001; I = 0.9994576682225345; J = 0.9994576682225345; O = 0.390990551238643; P = 0.390990551238643; Q = 2216.232226424801; R = 19946.0900378232; S = 0.3907785046380175; T = 0.3907785046380175 } },
    @{ Row = 8; Cells = @{ D = "MuSCs"; G = 214.0283; H = 642.0849000000001; I = 0.9994576682225345; J = 0.9994576682225345; M = 12.706793; N = 38.120379; O = 0.4797976910104138; P = 0.4797976910104139; Q = 2719.6133042419; R = 24476.5197381771; S = 0.4795374814758243; T = 0.4795374814758243 } },
    @{ Row = 9; Cells = @{ G = 214.0283; H = 642.0849000000001; I = 0.9994576682225345; J = 0.9994576682225345; M = 0.3439043333333334; N = 1.031713; O = 0.01298553498603535; P = 0.01298553498603535; Q = 73.60525982596668; R = 662.4473384337001; S = 0.01297849251776503; T = 0.01297849251776503 } }
)

foreach ($rowUpdate in $rowUpdates) {
    foreach ($col in $rowUpdate.Cells.Keys) {
        $addr = "$col$($rowUpdate.Row)"
        $ws.Range($addr).Value = $rowUpdate.Cells[$col]
    }
}
